# Daily roll-forward of the BP "terminal gate pricing" table.
#
# Each state section lists two (or more) dated rows per terminal: the
# most-recent effective date and the one before it. A new day's prices
# arrived, so for every terminal:
#   - the row that used to hold "yesterday"'s prices is dropped,
#   - the row that used to hold "today"'s prices becomes "yesterday"
#     (date shifts back one row, figures carried through unchanged), and
#   - the "today" row gets the freshly published prices (date advances
#     by one day).
#
# Column layout per data row: A = Effective Date (serial), B = (blank),
# C = Terminal name, D = Diesel, E = ULP, F = PULP, G = e10 (omitted
# where the terminal has no e10 price).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row=8; A=46044; D=157.58; E=150.74; F=160.74; G=150.75 },
    @{ Row=9; A=46044; D=157.58; E=150.74; F=160.74; G=150.75 },
    @{ Row=10; A=46044; D=158.36; E=152.51; F=162.51; G=152.91 },
    @{ Row=11; A=46043; D=157.16; E=150.8; F=160.8; G=150.82 },
    @{ Row=12; A=46043; D=157.16; E=150.8; F=160.8; G=150.82 },
    @{ Row=13; A=46043; D=158.1; E=152.73; F=162.73; G=153.13 },
    @{ Row=17; A=46044; D=161.74; E=155; F=165 },
    @{ Row=18; A=46043; D=161.51; E=155.25; F=165.25 },
    @{ Row=22; A=46044; D=158.65; E=152.31; F=161.91; G=153.38 },
    @{ Row=23; A=46044; D=163.14; E=157.84; F=167.84 },
    @{ Row=24; A=46044; D=163.27; E=158.56; F=168.56 },
    @{ Row=25; A=46044; D=163.25; E=158.11; F=168.11; G=158.23 },
    @{ Row=26; A=46044; D=162.83; E=159.72; F=169.72 },
    @{ Row=27; A=46043; D=158.35; E=152.49; F=162.09; G=153.56 },
    @{ Row=28; A=46043; D=162.87; E=158.06; F=168.06 },
    @{ Row=29; A=46043; D=163; E=158.77; F=168.77 },
    @{ Row=30; A=46043; D=162.99; E=158.31; F=168.31; G=158.44 },
    @{ Row=31; A=46043; D=162.58; E=159.93; F=169.93 },
    @{ Row=35; A=46044; D=157.24; E=149.5; F=158.5 },
    @{ Row=36; A=46043; D=156.98; E=149.72; F=158.72 },
    @{ Row=40; A=46044; D=162.8; E=157.8; F=167.8 },
    @{ Row=41; A=46044; D=162.52; E=158.22; F=168.22 },
    @{ Row=42; A=46043; D=162.52; E=158.05; F=168.05 },
    @{ Row=43; A=46043; D=162.23; E=158.47; F=168.47 },
    @{ Row=47; A=46044; D=156.76; E=151.16; F=161.16 },
    @{ Row=48; A=46044; D=156.35; E=151.07; F=161.07 },
    @{ Row=49; A=46043; D=156.33; E=151.27; F=161.27 },
    @{ Row=50; A=46043; D=155.93; E=151.18; F=161.18 },
    @{ Row=54; A=46044; D=171.8; E=165.87; F=175.87 },
    @{ Row=55; A=46044; D=164.61; E=163.7; F=173.7 },
    @{ Row=56; A=46044; D=161.32 },
    @{ Row=57; A=46044; D=161.64; E=158.12 },
    @{ Row=58; A=46044; D=157.41; E=154.02; F=164.02 },
    @{ Row=59; A=46044; D=164.05; E=163.87 },
    @{ Row=60; A=46043; D=171.55; E=165.97; F=175.97 },
    @{ Row=61; A=46043; D=164.3; E=163.95; F=173.95 },
    @{ Row=62; A=46043; D=161.04 },
    @{ Row=63; A=46043; D=161.42; E=158.37 },
    @{ Row=64; A=46043; D=157.18; E=154.27; F=164.27 },
    @{ Row=65; A=46043; D=163.82; E=164.02 }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.A
    if ($u.ContainsKey('D')) { $ws.Range("D$r").Value = $u.D }
    if ($u.ContainsKey('E')) { $ws.Range("E$r").Value = $u.E }
    if ($u.ContainsKey('F')) { $ws.Range("F$r").Value = $u.F }
    if ($u.ContainsKey('G')) { $ws.Range("G$r").Value = $u.G }
}
